# Adds core dry-mass data (CoreSleeveMass_g, SoilDryMass_g) to the
# "treatments.csv" sheet, matching the commit "Added core dry mass data".
#
# The original sheet had 30 "AL N" core rows (rows 5-34) below two
# "Ambient4"/"Ambient22" control rows (rows 3-4). The edit sorts the
# existing A:D block ascending by Core (column A) -- which naturally pushes
# the two Ambient rows to the bottom, rows 33-34 -- and then records the two
# new measurement columns for each "AL N" core at its now-sorted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers.
$ws.Range("E2").Value = "CoreSleeveMass_g"
$ws.Range("F2").Value = "SoilDryMass_g"

# Sort the existing table (A3:D34) ascending by Core (column A).
$sortRange = $ws.Range("A3:D34")
$keyRange = $ws.Range("A3:A34")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$ws.Sort.Apply()

# New measurements keyed by each Core's row position after the sort above
# (rows 33-34, the Ambient controls, get no new data).
$rowData = @{
    3  = @(148.9, 1759.312)     # AL 1
    4  = @(151.9, 832.484)      # AL 10
    5  = @(151.6, 1799.938)     # AL 12
    6  = @(150.9, 1229.824)     # AL 13
    7  = @(151.1, 1666.933)     # AL 14
    8  = @(149.3, 1421.643)     # AL 15
    9  = @(149.0, 1845.275)     # AL 16
    10 = @(147.7, 1495.182)     # AL 18
    11 = @(134.7, 1335.673)     # AL 19
    12 = @(154.8, 1964.05)      # AL 2
    13 = @(128.8, 1600.397)     # AL 20
    14 = @(148.4, 1907.169)     # AL 21
    15 = @(146.8, 1461.686)     # AL 22
    16 = @(129.9, 1708.298)     # AL 23
    17 = @(132.3, 1727.016)     # AL 24
    18 = @(149.6, 1162.577)     # AL 25
    19 = @(129.4, 1726.425)     # AL 26
    20 = @(148.2, 1140.263)     # AL 27
    21 = @(147.2, 1950.255)     # AL 28
    22 = @(154.9, 1351.708)     # AL 29
    23 = @(147.7, 1651.311)     # AL 34
    24 = @(155.4, 1486.275)     # AL 36
    25 = @(148.9, 1448.224)     # AL 37
    26 = @(149.4, 1567.252)     # AL 38
    27 = @(150.9, 1815.444)     # AL 39
    28 = @(152.1, 1414.568)     # AL 5
    29 = @(148.3, 1470.014)     # AL 6
    30 = @(133.5, 1640.926)     # AL 7
    31 = @(149.1, 735.107)      # AL 8
    32 = @(148.6, 1385.091)     # AL 9
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Cells.Item($r, 5).Value = $vals[0]
    $ws.Cells.Item($r, 6).Value = $vals[1]
}

# Selection marker matches the post-edit state recorded in the sheet.
$ws.Range("F2").Select()
